$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain plain text even when
# the new value looks like a number (e.g. "0.5303"), so Excel doesn't
# silently convert it to a numeric cell.
function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.167.30"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.653.12"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.47%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -0.20%  "

# Row 6 - XRP
Set-TextValue "D6" "0.5303"
$ws.Range("E6").Value = "  +0.61%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.39%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.2621"
$ws.Range("E8").Value = "  +0.10%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06333"
$ws.Range("E9").Value = "  +0.72%  "

# Row 10 - Solana
Set-TextValue "D10" "20.39"
$ws.Range("E10").Value = "  +0.69%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07798"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12 - Polkadot
Set-TextValue "D12" "4.518"
$ws.Range("E12").Value = "  +1.09%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.632.03"
$ws.Range("E13").Value = "  -2.16%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.881.11"
$ws.Range("E14").Value = "  +0.40%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.5491"
$ws.Range("E15").Value = "  +0.79%  "

# Row 16 - ShibaInu : no change

# Row 17 - Litecoin
Set-TextValue "D17" "65.45"
$ws.Range("E17").Value = "  +0.79%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "26.137.79"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19 - Dai : no change

# Row 20 - Uniswap
Set-TextValue "D20" "4.590"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21 - BitcoinCash
$ws.Range("E21").Value = "  -0.57%  "

# Row 22 - Avalanche
Set-TextValue "D22" "10.08"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.003"
$ws.Range("E23").Value = "  +0.44%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.45%  "

# Row 25 - Monero
Set-TextValue "D25" "145.17"
$ws.Range("E25").Value = "  +4.16%  "

# Row 26 - Stellar
Set-TextValue "D26" "0.1221"
$ws.Range("E26").Value = "  -1.64%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.89%  "

# Row 28 - EthereumClassic : no change

# Row 29 - Toncoin
Set-TextValue "D29" "1.476"
$ws.Range("E29").Value = "  +4.09%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.05731"
$ws.Range("E30").Value = "  -3.48%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.18%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "3.549"
$ws.Range("E32").Value = "  +1.70%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.261"
$ws.Range("E33").Value = "  +0.60%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.587"
$ws.Range("E34").Value = "  +3.36%  "

# Row 35 - MXToken
$ws.Range("E35").Value = "  +2.12%  "

# Row 36 - HuobiToken
Set-TextValue "D36" "2.422"
$ws.Range("E36").Value = "  +0.46%  "

# Row 37 - ARBITRUM
Set-TextValue "D37" "0.9482"
$ws.Range("E37").Value = "  +0.66%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.5747"
$ws.Range("E38").Value = "  +1.59%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.01602"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40 - TrustWalletToken
Set-TextValue "D40" "0.8488"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.779"
$ws.Range("E41").Value = "  -1.31%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.46%  "

# Row 43 - was Maker, now Quant
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "103.88"
$ws.Range("E43").Value = "  +2.95%  "

# Row 44 - was Quant, now Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.038.75"
$ws.Range("E44").Value = "  +3.19%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.793.97"
$ws.Range("E45").Value = "  +0.29%  "

# Row 46 - Aave
Set-TextValue "D46" "56.84"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  -1.35%  "

# Row 48 - Frax
Set-TextValue "D48" "1.007"
$ws.Range("E48").Value = "  +0.05%  "

# Row 49 - Mantle
Set-TextValue "D49" "0.4356"
$ws.Range("E49").Value = "  +1.72%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.874"
$ws.Range("E50").Value = "  +0.55%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.05154"
$ws.Range("E51").Value = "  +0.11%  "
